$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.878.28'
$ws.Range("E2").Value = '  +2.32%  '
$ws.Range("D3").Value = '3.931.87'
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '484.54'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("E7").Value = '  -1.69%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.723'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.169'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000354'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +12.00%  '
$ws.Range("E12").Value = '  -1.87%  '
$ws.Range("D13").Value = '4.571.33'
$ws.Range("E13").Value = '  +0.82%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.49'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.94%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.937.68'
$ws.Range("E15").Value = '  +1.58%  '
$ws.Range("B16").Value = 'Uniswap'
$ws.Range("C16").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.59'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.61%  '
$ws.Range("E17").Value = '  -0.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.69'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.07%  '
$ws.Range("E19").Value = '  -2.92%  '
$ws.Range("D20").Value = '69.003.89'
$ws.Range("E20").Value = '  +2.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '433.75'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.36'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.73%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.61'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.63'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.61'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +16.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.57'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '38.27'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.65%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.90'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '712.25'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.30'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.41%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.85'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.23%  '
$ws.Range("D34").Value = '0.0₃0917'
$ws.Range("E34").Value = '  +32.59%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '41.39'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.59'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.73%  '
$ws.Range("E37").Value = '  -7.39%  '
$ws.Range("E38").Value = '  +2.85%  '
$ws.Range("E39").Value = '  -0.14%  '
$ws.Range("E40").Value = '  -1.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.06'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +9.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.75'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.98'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.340'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.141'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("E46").Value = '  +0.22%  '
$ws.Range("E47").Value = '  -1.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '148.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.70%  '
$ws.Range("E50").Value = '  -4.53%  '
$ws.Range("E51").Value = '  -1.89%  '
